$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-28: price/volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.995.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.658.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.688.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.158.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.998.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.674.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "350.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.788.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +1.20%  "

# --- Rows 29-30: rank swap (InternetComputer now outranks PEPE) ---
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.56%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0810"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.77%  "

# --- Rows 31-49: price/volume updates ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.14%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.848"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "279.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0982"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0534"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +1.36%  "

# --- Rows 50-51: rank swap (Maker now outranks WhiteBITCoin) ---
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.002.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
